# Daily attendance processing - 2025-12-12 07:32:18
# Rotate the "Recorded By" (column G) list for every data row: move the
# first comma-separated name/email to the end of the list. Single-value
# cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Header is row 1 ("Recorded By"); data starts on row 2.
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $parts = $value -split ", "

    if ($parts.Length -gt 1) {
        $rest = $parts[1..($parts.Length - 1)]
        $rotated = $rest + $parts[0]
        $newValue = $rotated -join ", "
        $cell.Value = $newValue
    }
}
